# Applies the cell-value corrections captured in the target diff for the
# "CryCompanywiseStockReport" sheet: per-item Quantity (F) / Value (G) updates
# -- including a few same-item batch (B) row reorderings where Rate/MRP (D/E)
# travel with their row -- plus the resulting Sub Total / Company Total /
# Grand Total (B column) rollups.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 78
$ws.Range("G9").Value = 2306.46
$ws.Range("B10").Value = 36451.83
$ws.Range("F72").Value = 23
$ws.Range("G72").Value = 4811.14
$ws.Range("F90").Value = 69
$ws.Range("G90").Value = 9311.549999999999
$ws.Range("F109").Value = 105
$ws.Range("G109").Value = 13174.35
$ws.Range("B114").Value = 234984.85
$ws.Range("B146").Value = 53925
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 66.44
$ws.Range("B147").Value = 64350
$ws.Range("E147").Value = 70.63
$ws.Range("F147").Value = 2
$ws.Range("G147").Value = 132.88
$ws.Range("B148").Value = 57756
$ws.Range("E148").Value = 79.37
$ws.Range("F148").Value = -100
$ws.Range("G148").Value = -6644
$ws.Range("F150").Value = 265
$ws.Range("G150").Value = 25654.65
$ws.Range("B152").Value = 20133.41
$ws.Range("B163").Value = 64329
$ws.Range("E163").Value = 128.32
$ws.Range("F163").Value = 3
$ws.Range("G163").Value = 362.07
$ws.Range("B164").Value = 57552
$ws.Range("E164").Value = 136.86
$ws.Range("F164").Value = -5
$ws.Range("G164").Value = -603.45
$ws.Range("F186").Value = 1783
$ws.Range("G186").Value = 15066.35
$ws.Range("F187").Value = 975
$ws.Range("G187").Value = 7790.25
$ws.Range("B189").Value = 27739.62
$ws.Range("F193").Value = 280
$ws.Range("G193").Value = 18144
$ws.Range("B200").Value = 42111.1
$ws.Range("F255").Value = 28
$ws.Range("G255").Value = 897.96
$ws.Range("B274").Value = 67024.13
$ws.Range("B292").Value = 63520
$ws.Range("E292").Value = 153.4
$ws.Range("F292").Value = 69
$ws.Range("G292").Value = 9955.32
$ws.Range("B293").Value = 55373
$ws.Range("E293").Value = 163.62
$ws.Range("F293").Value = -94
$ws.Range("G293").Value = -13562.32
$ws.Range("B295").Value = 63571
$ws.Range("F295").Value = 0
$ws.Range("G295").Value = 0
$ws.Range("B296").Value = 63531
$ws.Range("F296").Value = 80
$ws.Range("G296").Value = 11478.4
$ws.Range("F304").Value = 15
$ws.Range("G304").Value = 4541.85
$ws.Range("F306").Value = 4
$ws.Range("G306").Value = 460.2
$ws.Range("F321").Value = 39
$ws.Range("G321").Value = 2289.3
$ws.Range("F328").Value = 221
$ws.Range("G328").Value = 4647.63
$ws.Range("F330").Value = 2
$ws.Range("G330").Value = 1051.9
$ws.Range("F334").Value = 341
$ws.Range("G334").Value = 51548.97
$ws.Range("B339").Value = 259354.29
$ws.Range("F373").Value = 75
$ws.Range("G373").Value = 2094.75
$ws.Range("F382").Value = 139
$ws.Range("G382").Value = 5978.39
$ws.Range("B395").Value = 226296.27
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 42
$ws.Range("G420").Value = 4432.68
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2
$ws.Range("F427").Value = 92
$ws.Range("G427").Value = 3423.32
$ws.Range("B430").Value = 38295.84
$ws.Range("F439").Value = 93
$ws.Range("G439").Value = 2204.1
$ws.Range("F442").Value = 28
$ws.Range("G442").Value = 1426.6
$ws.Range("F444").Value = 52
$ws.Range("G444").Value = 3742.44
$ws.Range("B448").Value = 35884.69
$ws.Range("B467").Value = 65068
$ws.Range("E467").Value = 13.97
$ws.Range("F467").Value = 63
$ws.Range("G467").Value = 828.45
$ws.Range("B468").Value = 53602
$ws.Range("E468").Value = 15.69
$ws.Range("F468").Value = -231
$ws.Range("G468").Value = -3037.65
$ws.Range("F480").Value = 123
$ws.Range("G480").Value = 1995.06
$ws.Range("F484").Value = 387
$ws.Range("G484").Value = 2511.63
$ws.Range("B487").Value = 45702
$ws.Range("E487").Value = 31.43
$ws.Range("F487").Value = -215
$ws.Range("G487").Value = -5654.5
$ws.Range("B488").Value = 64919
$ws.Range("E488").Value = 27.97
$ws.Range("F488").Value = 61
$ws.Range("G488").Value = 1604.3
$ws.Range("F490").Value = 189
$ws.Range("G490").Value = 2783.97
$ws.Range("B492").Value = -13575.79
$ws.Range("F498").Value = 89
$ws.Range("G498").Value = 5486.85
$ws.Range("B508").Value = 7112.78
$ws.Range("F515").Value = 3
$ws.Range("G515").Value = 110.31
$ws.Range("B528").Value = 16738.77
$ws.Range("B568").Value = 64810
$ws.Range("E568").Value = 291.22
$ws.Range("F568").Value = 5
$ws.Range("G568").Value = 1369.6
$ws.Range("B569").Value = 53319
$ws.Range("E569").Value = 310.64
$ws.Range("F569").Value = -6
$ws.Range("G569").Value = -1643.52
$ws.Range("F581").Value = 55
$ws.Range("G581").Value = 3498.55
$ws.Range("B586").Value = 19297.47
$ws.Range("B595").Value = 60031
$ws.Range("E595").Value = 111.69
$ws.Range("F595").Value = -5
$ws.Range("G595").Value = -492.5
$ws.Range("B596").Value = 64836
$ws.Range("E596").Value = 104.71
$ws.Range("F596").Value = 0
$ws.Range("G596").Value = 0
$ws.Range("F615").Value = 48
$ws.Range("G615").Value = 8422.559999999999
$ws.Range("B618").Value = 37306.62
$ws.Range("F620").Value = 22
$ws.Range("G620").Value = 2872.1
$ws.Range("F622").Value = 77
$ws.Range("G622").Value = 10052.35
$ws.Range("B625").Value = 18141.43
$ws.Range("F654").Value = 305
$ws.Range("G654").Value = 24515.9
$ws.Range("B655").Value = 32519.28
$ws.Range("F707").Value = 49
$ws.Range("G707").Value = 4987.71
$ws.Range("F712").Value = 43
$ws.Range("G712").Value = 1603.04
$ws.Range("F717").Value = 187
$ws.Range("G717").Value = 25246.87
$ws.Range("F719").Value = 100
$ws.Range("G719").Value = 12071
$ws.Range("B720").Value = 68723.23
$ws.Range("F741").Value = 122
$ws.Range("G741").Value = 5244.78
$ws.Range("F742").Value = 123
$ws.Range("G742").Value = 6136.47
$ws.Range("B747").Value = 55046.36
$ws.Range("F772").Value = 2705
$ws.Range("G772").Value = 441212.55
$ws.Range("B779").Value = 670287.33
$ws.Range("F782").Value = 49
$ws.Range("G782").Value = 7154.49
$ws.Range("B796").Value = 64338.35
$ws.Range("B797").Value = 2507804.15
$ws.Range("B798").Value = 2507804.15
